$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A156").Value = "IMX-USD"
$ws.Range("A157").Value = "GRT-USD"
